# El proyecto de horario queda terminado: se reparte el bloque del
# mediodia (antes 12:00-16:00 en tramos gruesos) en tramos mas finos
# de Estudio/Descanso/Ocio, y se anaden dos filas nuevas al final
# para "Familia y descanso" (20:00-21:00 y 21:00-22:00 pasan a ser
# franjas libres, y las dos ultimas horas del dia se documentan en
# las filas 20 y 21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fila de cabecera (sin cambios)
$header = @('Horas', 'Lunes', 'Martes', 'Miercoles', 'Jueves', 'Viernes', 'Sabado', 'Domingo')
for ($col = 0; $col -lt $header.Count; $col++) {
    $ws.Cells.Item(1, 1 + $col).Value = $header[$col]
}

# Tabla completa del horario: hora + actividad por dia (Lunes..Sabado).
# Una cadena vacia significa "sin actividad / celda en blanco".
$schedule = @(
    @{ Hora = '8:00 - 9:00';   Dias = @('', '', '', '', '', '') },
    @{ Hora = '9:00 - 10:00';  Dias = @('Clases', 'Clases', 'Clases', 'Clases', 'Clases', 'Clases') },
    @{ Hora = '10:00 - 11:00'; Dias = @('', '', '', '', '', '') },
    @{ Hora = '11:00 - 12:00'; Dias = @('', '', '', '', '', '') },
    @{ Hora = '12:00 - 12:30'; Dias = @('Clases', 'Clases', 'Clases', 'Clases', 'Clases', 'Clases') },
    @{ Hora = '12:30 - 13:00'; Dias = @('', '', '', '', 'Ocio', 'Ocio') },
    @{ Hora = '13:00 - 13:30'; Dias = @('Estudio', 'Estudio', '', '', '', '') },
    @{ Hora = '13:30 - 13:42'; Dias = @('Descanso', 'Descanso', '', '', '', '') },
    @{ Hora = '13:42 - 14:00'; Dias = @('Estudio', 'Estudio', '', '', '', '') },
    @{ Hora = '14:00 - 14:30'; Dias = @('Descanso', 'Descanso', '', '', 'Ocio', 'Ocio') },
    @{ Hora = '14:30 - 14:42'; Dias = @('Estudio', 'Estudio', '', '', '', '') },
    @{ Hora = '14:42 - 15:00'; Dias = @('', '', '', '', '', '') },
    @{ Hora = '15:00 - 16:00'; Dias = @('', '', '', '', '', '') },
    @{ Hora = '16:00 - 17:00'; Dias = @('', '', '', '', '', '') },
    @{ Hora = '17:00 - 18:00'; Dias = @('', '', '', '', '', '') },
    @{ Hora = '18:00 - 19:00'; Dias = @('', '', '', '', '', '') },
    @{ Hora = '19:00 - 20:00'; Dias = @('', '', '', '', '', '') },
    @{ Hora = '20:00 - 21:00'; Dias = @('', '', '', '', '', '') },
    @{ Hora = '21:00 - 22:00'; Dias = @('Familia y descanso', 'Familia y descanso', 'Familia y descanso', 'Familia y descanso', 'Familia y descanso', 'Familia y descanso') },
    @{ Hora = '22:00 - 23:00'; Dias = @('Familia y descanso', 'Familia y descanso', 'Familia y descanso', 'Familia y descanso', 'Familia y descanso', 'Familia y descanso') }
)

$row = 2
foreach ($slot in $schedule) {
    $ws.Cells.Item($row, 1).Value = $slot.Hora
    for ($col = 0; $col -lt $slot.Dias.Count; $col++) {
        $ws.Cells.Item($row, 2 + $col).Value = $slot.Dias[$col]
    }
    # Columna H (Domingo) siempre vacia en este horario
    $ws.Cells.Item($row, 8).Value = ''
    $row++
}
